$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) is treated as text so values like "1.000" or "0.8230" are preserved exactly
$ws.Range("D2:D51").NumberFormat = "@"

$data = @(
    ,@('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '30.182.03', '  +0.38%  ')
    ,@('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.911.21', '  -0.05%  ')
    ,@('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.000', '  -0.04%  ')
    ,@('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.8230', '  +4.25%  ')
    ,@('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '243.94', '  +0.38%  ')
    ,@('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.001', '  +0.01%  ')
    ,@('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3259', '  +2.75%  ')
    ,@('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '26.89', '  +1.98%  ')
    ,@('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07055', '  +1.79%  ')
    ,@('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.08103', '  +1.42%  ')
    ,@('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.7740', '  +3.45%  ')
    ,@('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.912.13', '  +0.02%  ')
    ,@('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.297', '  +1.11%  ')
    ,@('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '93.45', '  -0.06%  ')
    ,@('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '30.177.35', '  +0.32%  ')
    ,@('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '14.24', '  +1.32%  ')
    ,@('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.920', '  -0.37%  ')
    ,@('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '246.43', '  -0.39%  ')
    ,@('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007806', '  +0.11%  ')
    ,@('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.166.05', '  -0.09%  ')
    ,@('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.002', '  +0.11%  ')
    ,@('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  +0.01%  ')
    ,@('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.101', '  +2.53%  ')
    ,@('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1677', '  +21.01%  ')
    ,@('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.339', '  +0.15%  ')
    ,@('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '167.17', '  -1.71%  ')
    ,@('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '19.00', '  +0.27%  ')
    ,@('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.114', '  +3.26%  ')
    ,@('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.372', '  -0.25%  ')
    ,@('PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.528', '  +0.31%  ')
    ,@('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05849', '  +2.21%  ')
    ,@('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.312', '  -0.83%  ')
    ,@('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.108', '  -0.21%  ')
    ,@('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.276', '  +0.87%  ')
    ,@('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7361', '  -0.22%  ')
    ,@('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.694', '  -1.30%  ')
    ,@('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01927', '  -0.37%  ')
    ,@('MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.798', '  +0.12%  ')
    ,@('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4469', '  +0.40%  ')
    ,@('Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '73.51', '  +1.21%  ')
    ,@('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '5.975', '  -3.36%  ')
    ,@('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8533', '  +2.13%  ')
    ,@('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.915', '  +0.66%  ')
    ,@('PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.001', '  -0.03%  ')
    ,@('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '102.75', '  +2.13%  ')
    ,@('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.600', '  +0.22%  ')
    ,@('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.862', '  +0.29%  ')
    ,@('Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.006.33', '  +1.81%  ')
    ,@('RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '2.064.54', '  +0.14%  ')
    ,@('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.562', '  +3.95%  ')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}

Write-Host "Updated cryptos table rows 2-51"
